# Apply updated values to the TDG balance sheet (row 4: Inventory,
# row 15: Accounts Payable, row 23: Long Term Tax Liability (Deferred))
# for columns B through F.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 - Inventory
$ws.Range("B4").Value = 1284000000.0
$ws.Range("C4").Value = 1283000000.0
$ws.Range("D4").Value = 1344000000.0
$ws.Range("E4").Value = 1313000000.0
$ws.Range("F4").Value = 1294000000.0

# Row 15 - Accounts Payable
$ws.Range("B15").Value = 197000000.0
$ws.Range("C15").Value = 218000000.0
$ws.Range("D15").Value = 231000000.0
$ws.Range("E15").Value = 266000000.0
$ws.Range("F15").Value = 265000000.0

# Row 23 - Long Term Tax Liability (Deferred)
$ws.Range("B23").Value = 419000000.0
$ws.Range("C23").Value = 413000000.0
$ws.Range("D23").Value = 360000000.0
$ws.Range("E23").Value = 372000000.0
$ws.Range("F23").Value = 435000000.0
